$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data for rows 2-11 (A: id, B: name, C: department, D: reason, E: hours, F: date serial, G: salary)
$data = @(
    @(74531, "Gustavo Silveira", "TI", "Viagem de negócios", 1, 45082, 5734.1),
    @(10235, "Theo da Mata", "Engenharia", "Problemas pessoais", 3, 45080, 11726.9),
    @(72835, "Dr. Lucas Gabriel Costela", "P&D", "Viagem de negócios", 2, 45099, 11677.71),
    @(50082, "Agatha Novaes", "Jurídico", "Doença", 1, 45100, 10912.45),
    @(73151, "Sr. Marcos Vinicius Costela", "Vendas", "Problemas pessoais", 5, 45106, 9916.120000000001),
    @(64764, "Lorenzo da Cruz", "Recursos Humanos", "Outros", 3, 45101, 5335.08),
    @(83659, "Lorenzo Dias", "Jurídico", "Problemas pessoais", 4, 45097, 8541.6),
    @(84162, "Lorenzo Costela", "Engenharia", "Problemas pessoais", 1, 45099, 8665.25),
    @(42794, "Vitor Pinto", "P&D", "Consulta médica", 4, 45088, 11318.1),
    @(16411, "Rafaela Barbosa", "Recursos Humanos", "Problemas pessoais", 2, 45088, 8708.790000000001)
)

$rowIndex = 2
foreach ($rec in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $rec[0]
    $ws.Cells.Item($rowIndex, 2).Value = $rec[1]
    $ws.Cells.Item($rowIndex, 3).Value = $rec[2]
    $ws.Cells.Item($rowIndex, 4).Value = $rec[3]
    $ws.Cells.Item($rowIndex, 5).Value = $rec[4]
    $ws.Cells.Item($rowIndex, 6).Value = $rec[5]
    $ws.Cells.Item($rowIndex, 7).Value = $rec[6]
    $rowIndex++
}
